# Apply the gh-pages metadata refresh to ValueSet-tnm-stage-group-vs.xlsx
#  - rename the three "Include from SNOMED CT[ N]" sheets to "Include #N"
#  - bump the Metadata sheet's Date value
#  - insert a new "Jurisdiction" property row (empty value) right before
#    "Description", pushing Description/Purpose/Copyright/Immutable down
#    by one row

$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from SNOMED CT" worksheets -------------------
$wb.Worksheets.Item("Include from SNOMED CT").Name   = "Include #0"
$wb.Worksheets.Item("Include from SNOMED CT 2").Name = "Include #1"
$wb.Worksheets.Item("Include from SNOMED CT 3").Name = "Include #2"

# --- 2. Metadata sheet updates --------------------------------------------
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row above "Description" (row 11) for "Jurisdiction".
$ws.Rows.Item(11).Insert()

# Match the formatting of the surrounding data rows (border/alignment).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh the publication Date value.
$ws.Range("B8").Value = "2024-09-17T19:55:11+00:00"
